# Update the "Total" (column O) values on Sheet1 rows 4-36 with refreshed
# figures, and move the active selection to X12 (as left by the author
# after editing), mirroring the authoritative diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O4").Value  = -1998.5981221325214
$ws.Range("O5").Value  = 16435.781526383507
$ws.Range("O6").Value  = 20454.849140960243
$ws.Range("O7").Value  = 17038.94422737664
$ws.Range("O8").Value  = 13745.047311162791
$ws.Range("O9").Value  = 15132.06000946791
$ws.Range("O10").Value = 15606.987509797171
$ws.Range("O11").Value = 17016.485124051971
$ws.Range("O12").Value = 17591.634510507723
$ws.Range("O13").Value = 18498.030802662179
$ws.Range("O14").Value = 12827.807462571316
$ws.Range("O15").Value = 13790.437082697636
$ws.Range("O16").Value = 15298.703582549764
$ws.Range("O17").Value = 16564.725043283786
$ws.Range("O18").Value = 17848.879726987452
$ws.Range("O19").Value = 18757.600227044673
$ws.Range("O20").Value = 20290.664987979435
$ws.Range("O21").Value = 27634.185731017446
$ws.Range("O22").Value = 18079.404544884062
$ws.Range("O23").Value = 18649.553918816058
$ws.Range("O24").Value = 18956.452144169099
$ws.Range("O25").Value = 19240.652823209126
$ws.Range("O26").Value = 17867.330379917476
$ws.Range("O27").Value = 17166.372303582008
$ws.Range("O28").Value = 18060.034715440248
$ws.Range("O29").Value = 19741.439858855949
$ws.Range("O30").Value = 19565.728775564799
$ws.Range("O31").Value = 19658.773743151509
$ws.Range("O32").Value = 6505.617728611267
$ws.Range("O33").Value = 5197.5713810805883
$ws.Range("O34").Value = 4230.5051756939165
$ws.Range("O35").Value = 3259.3450200306347
$ws.Range("O36").Value = 1144.3268849167191

# Leave the selection where the author left it after making the edits.
$ws.Range("X12").Select()
